# Apply the NATMI TPM update to the Adm2-Ramp3 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 5, target cluster "Resolving-Mac") entirely.
$ws.Rows.Item(5).Delete()

# Row 4's target cluster label changes from "MuSCs" to "Inflammatory-Mac".
$ws.Range("D4").Value = "Inflammatory-Mac"

# --- Row 2 (target cluster: ECs) ---
$ws.Range("H2").Value = 0.5797639999999999
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5477486666666667
$ws.Range("N2").Value = 1.643246
$ws.Range("O2").Value = 0.430739219238867
$ws.Range("P2").Value = 0.430739219238867
$ws.Range("Q2").Value = 0.1058549859937778
$ws.Range("R2").Value = 0.9526948739439999
$ws.Range("S2").Value = 0.430739219238867
$ws.Range("T2").Value = 0.430739219238867

# --- Row 3 (target cluster: FAPs) ---
$ws.Range("H3").Value = 0.5797639999999999
$ws.Range("M3").Value = 0.483825
$ws.Range("O3").Value = 0.3804708535695413
$ws.Range("P3").Value = 0.3804708535695413
$ws.Range("S3").Value = 0.3804708535695413
$ws.Range("T3").Value = 0.3804708535695413

# --- Row 4 (target cluster: Inflammatory-Mac, formerly MuSCs) ---
$ws.Range("H4").Value = 0.5797639999999999
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2400743333333334
$ws.Range("N4").Value = 0.7202230000000001
$ws.Range("O4").Value = 0.1887899271915918
$ws.Range("P4").Value = 0.1887899271915918
$ws.Range("Q4").Value = 0.04639548526355556
$ws.Range("R4").Value = 0.417559367372
$ws.Range("S4").Value = 0.1887899271915918
$ws.Range("T4").Value = 0.1887899271915918
